# sf9.xlsx template update:
#  - Fill in dummy quarterly grades for the student grades table (BACK sheet)
#  - Wire up the MAPEH average row (17) and its sub-rows (18-21) with formulas
#  - Wire up the General Average row (22) with its formula + PASSED/FAILED remark
#  - Update selection to reflect where the user ended up after editing

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BACK")

# --- Quarterly grade entries (columns N,O,P,Q) + Final rating (R) for the
#     single-quarter subjects that already had quarter grades filled in ---
$ws.Range("R7").Value = 90
$ws.Range("R8").Value = 90
$ws.Range("R9").Value = 90
$ws.Range("R11").Value = 90
$ws.Range("R12").Value = 90
$ws.Range("R13").Value = 90
$ws.Range("R15").Value = 90

# --- MAPEH (row 17) averages its 4 sub-components (rows 18-21) ---
$ws.Range("N17").Formula = "=SUM(N18:N21)/4"
$ws.Range("O17").Formula = "=SUM(O18:O21)/4"
$ws.Range("P17").Formula = "=SUM(P18:P21)/4"
$ws.Range("Q17").Formula = "=SUM(Q18:Q21)/4"
$ws.Range("R17").Formula = "=SUM(N17:Q17)/4"
$ws.Range("S17").Formula = '=IF(R17>=74, "PASSED", "FAILED")'

# --- MAPEH sub-component dummy grades ---
$ws.Range("N18:Q18").Value = 75
$ws.Range("N19:Q19").Value = 85
$ws.Range("N20:Q20").Value = 75
$ws.Range("N21:Q21").Value = 75

$ws.Range("R18").Formula = "=SUM(N18:Q18)/4"
$ws.Range("R19").Formula = "=SUM(N19:Q19)/4"
$ws.Range("R20").Formula = "=SUM(N20:Q20)/4"
$ws.Range("R21").Formula = "=SUM(N21:Q21)/4"

$ws.Range("S18").Formula = '=IF(R18>75, "FAILED", "PASSED")'
$ws.Range("S19").Formula = '=IF(R19>75, "FAILED", "PASSED")'
$ws.Range("S20").Formula = '=IF(R20>75, "FAILED", "PASSED")'
$ws.Range("S21").Formula = '=IF(R21>75, "FAILED", "PASSED")'

# --- General Average (row 22) ---
$ws.Range("R22").Formula = "=(R7+R8+R9+R11+R12+R13+R15)/7"
$ws.Range("S22").Formula = '=IF(R22>=74, "PASSED", "FAILED")'

# --- Formatting clean-up: match the Final-Rating/Remarks column styling that
#     is already used by the neighbouring cells in these rows ---
$ws.Range("S13").Copy()
$ws.Range("R13").PasteSpecial(-4122)
$ws.Range("S14").Copy()
$ws.Range("R14").PasteSpecial(-4122)
$ws.Range("S15").Copy()
$ws.Range("R15").PasteSpecial(-4122)
$ws.Range("S16").Copy()
$ws.Range("R16").PasteSpecial(-4122)

$ws.Range("R17").Copy()
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("R18:S21").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Leave the selection where the user's edits ended up ---
$ws.Range("S18").Select()
